$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "245.97"
$ws.Range("G2").Value = "12"
$ws.Range("G3").Value = "12"
$ws.Range("D4").Value = "5.426"
$ws.Range("G4").Value = "12"
$ws.Range("D5").Value = "0.05841"
$ws.Range("G5").Value = "12"
$ws.Range("D6").Value = "3.376"
$ws.Range("G6").Value = "12"
$ws.Range("D7").Value = "6.339"
$ws.Range("G7").Value = "12"
$ws.Range("D8").Value = "0.8080"
$ws.Range("G8").Value = "12"
$ws.Range("D9").Value = "0.9726"
$ws.Range("G9").Value = "12"
$ws.Range("D10").Value = "0.1426"
$ws.Range("G10").Value = "12"
$ws.Range("D11").Value = "0.07506"
$ws.Range("G11").Value = "12"
$ws.Range("D12").Value = "0.03372"
$ws.Range("G12").Value = "12"
$ws.Range("D13").Value = "0.02999"
$ws.Range("G13").Value = "12"
$ws.Range("D14").Value = "4.151"
$ws.Range("G14").Value = "12"
$ws.Range("G15").Value = "12"
$ws.Range("D16").Value = "0.001586"
$ws.Range("G16").Value = "12"
$ws.Range("D17").Value = "0.04808"
$ws.Range("G17").Value = "12"
$ws.Range("D18").Value = "0.0005891"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("G18").Value = "12"
$ws.Range("D19").Value = "0.006107"
$ws.Range("G19").Value = "12"
$ws.Range("G20").Value = "12"
$ws.Range("D21").Value = "0.0009990"
$ws.Range("G21").Value = "12"
$ws.Range("G22").Value = "12"
$ws.Range("D23").Value = "3.699"
$ws.Range("G23").Value = "12"
$ws.Range("D24").Value = "2.222"
$ws.Range("G24").Value = "12"
$ws.Range("G25").Value = "12"
$ws.Range("G26").Value = "12"
$ws.Range("D27").Value = "0.0001291"
$ws.Range("E27").Value = "26UpBotsUBXTWorstin24h"
$ws.Range("G27").Value = "12"
$ws.Range("G28").Value = "12"
$ws.Range("G29").Value = "12"
$ws.Range("G30").Value = "12"
$ws.Range("G31").Value = "12"
$ws.Range("G32").Value = "12"
$ws.Range("G33").Value = "12"
$ws.Range("G34").Value = "12"
$ws.Range("G35").Value = "12"
$ws.Range("G36").Value = "12"
$ws.Range("G37").Value = "12"
$ws.Range("G38").Value = "12"
$ws.Range("G39").Value = "12"
$ws.Range("D40").Value = "0.03859"
$ws.Range("G40").Value = "12"
$ws.Range("G41").Value = "12"
$ws.Range("D42").Value = "0.002441"
$ws.Range("G42").Value = "12"
$ws.Range("D43").Value = "0.003032"
$ws.Range("G43").Value = "12"
$ws.Range("D44").Value = "0.006686"
$ws.Range("G44").Value = "12"
$ws.Range("D45").Value = "0.00005617"
$ws.Range("G45").Value = "12"
$ws.Range("G46").Value = "12"
$ws.Range("D47").Value = "0.3901"
$ws.Range("G47").Value = "12"
$ws.Range("D48").Value = "0.1429"
$ws.Range("G48").Value = "12"
$ws.Range("G49").Value = "12"
$ws.Range("G50").Value = "12"
$ws.Range("G51").Value = "12"
